# Update the NATMI Slit3-Robo2 LR-pairs sheet with newly computed TPM-based
# statistics. The sending/target cluster combinations now form the full 4x4
# cross product of ECs / FAPs / MuSCs / Resolving-Mac (16 data rows, rows 2-17)
# instead of the previous partial 12-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Slit3"
$ws.Cells.Item(2,3).Value = "Robo2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.6868273333333333
$ws.Cells.Item(2,8).Value = 2.060482
$ws.Cells.Item(2,9).Value = 0.01130642661970366
$ws.Cells.Item(2,10).Value = 0.01130642661970366
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.9052683333333333
$ws.Cells.Item(2,14).Value = 2.715805
$ws.Cells.Item(2,15).Value = 0.5200140314301739
$ws.Cells.Item(2,16).Value = 0.5200140314301738
$ws.Cells.Item(2,17).Value = 0.6217630353344444
$ws.Cells.Item(2,18).Value = 5.59586731801
$ws.Cells.Item(2,19).Value = 0.005879500487581534
$ws.Cells.Item(2,20).Value = 0.005879500487581533

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Slit3"
$ws.Cells.Item(3,3).Value = "Robo2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.6868273333333333
$ws.Cells.Item(3,8).Value = 2.060482
$ws.Cells.Item(3,9).Value = 0.01130642661970366
$ws.Cells.Item(3,10).Value = 0.01130642661970366
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.7756663333333332
$ws.Cells.Item(3,14).Value = 2.326999
$ws.Cells.Item(3,15).Value = 0.4455666482402024
$ws.Cells.Item(3,16).Value = 0.4455666482402024
$ws.Cells.Item(3,17).Value = 0.5327488392797777
$ws.Cells.Item(3,18).Value = 4.794739553517999
$ws.Cells.Item(3,19).Value = 0.005037766612515162
$ws.Cells.Item(3,20).Value = 0.005037766612515161

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Slit3"
$ws.Cells.Item(4,3).Value = "Robo2"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.6868273333333333
$ws.Cells.Item(4,8).Value = 2.060482
$ws.Cells.Item(4,9).Value = 0.01130642661970366
$ws.Cells.Item(4,10).Value = 0.01130642661970366
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.04422333333333334
$ws.Cells.Item(4,14).Value = 0.13267
$ws.Cells.Item(4,15).Value = 0.02540324564902162
$ws.Cells.Item(4,16).Value = 0.02540324564902162
$ws.Cells.Item(4,17).Value = 0.03037379410444445
$ws.Cells.Item(4,18).Value = 0.27336414694
$ws.Cells.Item(4,19).Value = 0.0002872199328329693
$ws.Cells.Item(4,20).Value = 0.0002872199328329692

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Slit3"
$ws.Cells.Item(5,3).Value = "Robo2"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.6868273333333333
$ws.Cells.Item(5,8).Value = 2.060482
$ws.Cells.Item(5,9).Value = 0.01130642661970366
$ws.Cells.Item(5,10).Value = 0.01130642661970366
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.01569566666666667
$ws.Cells.Item(5,14).Value = 0.047087
$ws.Cells.Item(5,15).Value = 0.009016074680602103
$ws.Cells.Item(5,16).Value = 0.009016074680602101
$ws.Cells.Item(5,17).Value = 0.01078021288155556
$ws.Cells.Item(5,18).Value = 0.09702191593399999
$ws.Cells.Item(5,19).Value = 0.0001019395867739958
$ws.Cells.Item(5,20).Value = 0.0001019395867739958

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Slit3"
$ws.Cells.Item(6,3).Value = "Robo2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 53.540432
$ws.Cells.Item(6,8).Value = 160.621296
$ws.Cells.Item(6,9).Value = 0.8813728519762372
$ws.Cells.Item(6,10).Value = 0.881372851976237
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.9052683333333333
$ws.Cells.Item(6,14).Value = 2.715805
$ws.Cells.Item(6,15).Value = 0.5200140314301739
$ws.Cells.Item(6,16).Value = 0.5200140314301738
$ws.Cells.Item(6,17).Value = 48.46845764258667
$ws.Cells.Item(6,18).Value = 436.21611878328
$ws.Cells.Item(6,19).Value = 0.458326249949273
$ws.Cells.Item(6,20).Value = 0.4583262499492728

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Slit3"
$ws.Cells.Item(7,3).Value = "Robo2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 53.540432
$ws.Cells.Item(7,8).Value = 160.621296
$ws.Cells.Item(7,9).Value = 0.8813728519762372
$ws.Cells.Item(7,10).Value = 0.881372851976237
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7756663333333332
$ws.Cells.Item(7,14).Value = 2.326999
$ws.Cells.Item(7,15).Value = 0.4455666482402024
$ws.Cells.Item(7,16).Value = 0.4455666482402024
$ws.Cells.Item(7,17).Value = 41.52951057452266
$ws.Cells.Item(7,18).Value = 373.765595170704
$ws.Cells.Item(7,19).Value = 0.3927103475049601
$ws.Cells.Item(7,20).Value = 0.39271034750496

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Slit3"
$ws.Cells.Item(8,3).Value = "Robo2"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 53.540432
$ws.Cells.Item(8,8).Value = 160.621296
$ws.Cells.Item(8,9).Value = 0.8813728519762372
$ws.Cells.Item(8,10).Value = 0.881372851976237
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.04422333333333334
$ws.Cells.Item(8,14).Value = 0.13267
$ws.Cells.Item(8,15).Value = 0.02540324564902162
$ws.Cells.Item(8,16).Value = 0.02540324564902162
$ws.Cells.Item(8,17).Value = 2.367736371146667
$ws.Cells.Item(8,18).Value = 21.30962734032
$ws.Cells.Item(8,19).Value = 0.02238973106713113
$ws.Cells.Item(8,20).Value = 0.02238973106713112

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Slit3"
$ws.Cells.Item(9,3).Value = "Robo2"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 53.540432
$ws.Cells.Item(9,8).Value = 160.621296
$ws.Cells.Item(9,9).Value = 0.8813728519762372
$ws.Cells.Item(9,10).Value = 0.881372851976237
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.01569566666666667
$ws.Cells.Item(9,14).Value = 0.047087
$ws.Cells.Item(9,15).Value = 0.009016074680602103
$ws.Cells.Item(9,16).Value = 0.009016074680602101
$ws.Cells.Item(9,17).Value = 0.8403527738613333
$ws.Cells.Item(9,18).Value = 7.563174964751999
$ws.Cells.Item(9,19).Value = 0.007946523454873016
$ws.Cells.Item(9,20).Value = 0.007946523454873015

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Slit3"
$ws.Cells.Item(10,3).Value = "Robo2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 6.476716
$ws.Cells.Item(10,8).Value = 19.430148
$ws.Cells.Item(10,9).Value = 0.1066185206043934
$ws.Cells.Item(10,10).Value = 0.1066185206043934
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.9052683333333333
$ws.Cells.Item(10,14).Value = 2.715805
$ws.Cells.Item(10,15).Value = 0.5200140314301739
$ws.Cells.Item(10,16).Value = 0.5200140314301738
$ws.Cells.Item(10,17).Value = 5.863165898793333
$ws.Cells.Item(10,18).Value = 52.76849308914
$ws.Cells.Item(10,19).Value = 0.05544312672461169
$ws.Cells.Item(10,20).Value = 0.05544312672461169

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Slit3"
$ws.Cells.Item(11,3).Value = "Robo2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 6.476716
$ws.Cells.Item(11,8).Value = 19.430148
$ws.Cells.Item(11,9).Value = 0.1066185206043934
$ws.Cells.Item(11,10).Value = 0.1066185206043934
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.7756663333333332
$ws.Cells.Item(11,14).Value = 2.326999
$ws.Cells.Item(11,15).Value = 0.4455666482402024
$ws.Cells.Item(11,16).Value = 0.4455666482402024
$ws.Cells.Item(11,17).Value = 5.023770551761332
$ws.Cells.Item(11,18).Value = 45.213934965852
$ws.Cells.Item(11,19).Value = 0.04750565686602855
$ws.Cells.Item(11,20).Value = 0.04750565686602854

# Row 12
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Slit3"
$ws.Cells.Item(12,3).Value = "Robo2"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 6.476716
$ws.Cells.Item(12,8).Value = 19.430148
$ws.Cells.Item(12,9).Value = 0.1066185206043934
$ws.Cells.Item(12,10).Value = 0.1066185206043934
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.04422333333333334
$ws.Cells.Item(12,14).Value = 0.13267
$ws.Cells.Item(12,15).Value = 0.02540324564902162
$ws.Cells.Item(12,16).Value = 0.02540324564902162
$ws.Cells.Item(12,17).Value = 0.2864219705733334
$ws.Cells.Item(12,18).Value = 2.57779773516
$ws.Cells.Item(12,19).Value = 0.00270845646964868
$ws.Cells.Item(12,20).Value = 0.002708456469648679

# Row 13
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Slit3"
$ws.Cells.Item(13,3).Value = "Robo2"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 6.476716
$ws.Cells.Item(13,8).Value = 19.430148
$ws.Cells.Item(13,9).Value = 0.1066185206043934
$ws.Cells.Item(13,10).Value = 0.1066185206043934
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.01569566666666667
$ws.Cells.Item(13,14).Value = 0.047087
$ws.Cells.Item(13,15).Value = 0.009016074680602103
$ws.Cells.Item(13,16).Value = 0.009016074680602101
$ws.Cells.Item(13,17).Value = 0.1016563754306667
$ws.Cells.Item(13,18).Value = 0.9149073788759999
$ws.Cells.Item(13,19).Value = 0.0009612805441045254
$ws.Cells.Item(13,20).Value = 0.0009612805441045252

# Row 14
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Slit3"
$ws.Cells.Item(14,3).Value = "Robo2"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.04265633333333333
$ws.Cells.Item(14,8).Value = 0.127969
$ws.Cells.Item(14,9).Value = 0.0007022007996657373
$ws.Cells.Item(14,10).Value = 0.0007022007996657372
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.9052683333333333
$ws.Cells.Item(14,14).Value = 2.715805
$ws.Cells.Item(14,15).Value = 0.5200140314301739
$ws.Cells.Item(14,16).Value = 0.5200140314301738
$ws.Cells.Item(14,17).Value = 0.03861542778277777
$ws.Cells.Item(14,18).Value = 0.347538850045
$ws.Cells.Item(14,19).Value = 0.000365154268707672
$ws.Cells.Item(14,20).Value = 0.0003651542687076718

# Row 15
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Slit3"
$ws.Cells.Item(15,3).Value = "Robo2"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.04265633333333333
$ws.Cells.Item(15,8).Value = 0.127969
$ws.Cells.Item(15,9).Value = 0.0007022007996657373
$ws.Cells.Item(15,10).Value = 0.0007022007996657372
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.7756663333333332
$ws.Cells.Item(15,14).Value = 2.326999
$ws.Cells.Item(15,15).Value = 0.4455666482402024
$ws.Cells.Item(15,16).Value = 0.4455666482402024
$ws.Cells.Item(15,17).Value = 0.0330870816701111
$ws.Cells.Item(15,18).Value = 0.297783735031
$ws.Cells.Item(15,19).Value = 0.0003128772566986524
$ws.Cells.Item(15,20).Value = 0.0003128772566986523

# Row 16
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Slit3"
$ws.Cells.Item(16,3).Value = "Robo2"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.04265633333333333
$ws.Cells.Item(16,8).Value = 0.127969
$ws.Cells.Item(16,9).Value = 0.0007022007996657373
$ws.Cells.Item(16,10).Value = 0.0007022007996657372
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 0.6666666666666666
$ws.Cells.Item(16,13).Value = 0.04422333333333334
$ws.Cells.Item(16,14).Value = 0.13267
$ws.Cells.Item(16,15).Value = 0.02540324564902162
$ws.Cells.Item(16,16).Value = 0.02540324564902162
$ws.Cells.Item(16,17).Value = 0.001886405247777778
$ws.Cells.Item(16,18).Value = 0.01697764723
$ws.Cells.Item(16,19).Value = 0.00001783817940884815
$ws.Cells.Item(16,20).Value = 0.00001783817940884814

# Row 17
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Slit3"
$ws.Cells.Item(17,3).Value = "Robo2"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.04265633333333333
$ws.Cells.Item(17,8).Value = 0.127969
$ws.Cells.Item(17,9).Value = 0.0007022007996657373
$ws.Cells.Item(17,10).Value = 0.0007022007996657372
$ws.Cells.Item(17,11).Value = 1
$ws.Cells.Item(17,12).Value = 0.3333333333333333
$ws.Cells.Item(17,13).Value = 0.01569566666666667
$ws.Cells.Item(17,14).Value = 0.047087
$ws.Cells.Item(17,15).Value = 0.009016074680602103
$ws.Cells.Item(17,16).Value = 0.009016074680602101
$ws.Cells.Item(17,17).Value = 0.0006695195892222222
$ws.Cells.Item(17,18).Value = 0.006025676303
$ws.Cells.Item(17,19).Value = 0.000006331094850564804
$ws.Cells.Item(17,20).Value = 0.000006331094850564802
